$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.004866600036621
$ws.Range("B1").Value = 2.1470627784729
$ws.Range("C1").Value = 7.28476095199585
$ws.Range("D1").Value = 2.388560056686401
$ws.Range("E1").Value = 1.344383955001831
